# Update cryptocurrency price/volume data per the latest GitHub Actions refresh.
# Row 37/38 additionally swap (Monero <-> EthereumClassic changed rank order).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.225.60'
$ws.Range('E2').Value = '  -0.17%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.643.26'
$ws.Range('E3').Value = '  +0.11%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.54'
$ws.Range('E5').Value = '  -0.41%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.18'
$ws.Range('E6').Value = '  +2.03%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.542'
$ws.Range('E8').Value = '  -0.73%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.141'
$ws.Range('E9').Value = '  +2.27%  '

$ws.Range('E10').Value = '  -1.33%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.24'
$ws.Range('E11').Value = '  +0.56%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.351'
$ws.Range('E12').Value = '  +0.87%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.02'
$ws.Range('E13').Value = '  +0.53%  '

$ws.Range('E14').Value = '  +0.89%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.125.37'
$ws.Range('E15').Value = '  -0.10%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '68.155.34'
$ws.Range('E16').Value = '  -0.11%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.669.74'
$ws.Range('E17').Value = '  +1.33%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.37'
$ws.Range('E18').Value = '  -0.25%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '363.60'
$ws.Range('E19').Value = '  -0.26%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.32'
$ws.Range('E20').Value = '  -0.92%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.41'
$ws.Range('E21').Value = '  +3.58%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.79'
$ws.Range('E22').Value = '  -0.98%  '

$ws.Range('E23').Value = '  -1.87%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '75.18'
$ws.Range('E24').Value = '  +2.48%  '

$ws.Range('E25').Value = '  +0.17%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.70'
$ws.Range('E26').Value = '  -2.53%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.778.44'
$ws.Range('E27').Value = '  +0.18%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000104'
$ws.Range('E28').Value = '  -0.12%  '

$ws.Range('E29').Value = '  +0.19%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '556.52'
$ws.Range('E30').Value = '  -2.70%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.01'
$ws.Range('E31').Value = '  +0.58%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.40'
$ws.Range('E32').Value = '  -1.44%  '

$ws.Range('E33').Value = '  +0.28%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.128'
$ws.Range('E34').Value = '  -0.39%  '

$ws.Range('E35').Value = '  -0.04%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.55'
$ws.Range('E36').Value = '  +0.89%  '

$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.82'
$ws.Range('E37').Value = '  +3.17%  '

$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '160.04'
$ws.Range('E38').Value = '  -0.09%  '

$ws.Range('E39').Value = '  +0.80%  '

$ws.Range('E40').Value = '  -2.44%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.32'
$ws.Range('E41').Value = '  -0.54%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0₆0331'
$ws.Range('E42').Value = '  +4.04%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.59'
$ws.Range('E43').Value = '  -0.60%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '158.77'
$ws.Range('E45').Value = '  +1.26%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.74'
$ws.Range('E46').Value = '  +0.25%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '22.06'
$ws.Range('E47').Value = '  +1.06%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.68'
$ws.Range('E48').Value = '  -0.94%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0780'
$ws.Range('E49').Value = '  +0.22%  '

$ws.Range('E50').Value = '  -0.28%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.565'
$ws.Range('E51').Value = '  +0.68%  '
